$wb = $excel.ActiveWorkbook

function Set-Cells($ws, $row, $cols, $vals) {
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Cells.Item($row, $cols[$i]).Value = $vals[$i]
    }
}

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
Set-Cells $ws 5 @(8, 9, 10, 11, 12, 13, 14) @(180.15384, 86.666664, 390.5, 86.666664, 390.5, 28.333336, -620.5)  # HIJKLMN5
Set-Cells $ws 7 @(8, 10, 12, 14) @(2999.5, 2999, 2999, -3223)  # HJLN7
Set-Cells $ws 10 @(8, 9, 11, 13) @(500, 500, 500, -207)  # HIKM10
Set-Cells $ws 13 @(8, 10, 12, 14) @(10000, 10000, 10000, -10338)  # HJLN13
Set-Cells $ws 14 @(8, 10, 12, 14) @(2999.5, 2999, 2999, -3381)  # HJLN14
Set-Cells $ws 75 @(8, 10, 12, 14) @(33304.668, 40314, 40314, -42186)  # HJLN75
Set-Cells $ws 78 @(8, 10, 12, 14) @(33304.668, 40314, 120942, -130302)  # HJLN78
Set-Cells $ws 80 @(8, 9, 10, 11, 12, 13, 14) @(4905.0415, 238.83333, 18903.666, 716.49999, 56710.99800000001, 281.50001, -58706.99800000001)  # HIJKLMN80
Set-Cells $ws 83 @(8, 9, 10, 11, 12, 13, 14) @(4905.0415, 238.83333, 18903.666, 2149.49997, 170132.994, 2842.50003, -180116.994)  # HIJKLMN83
Set-Cells $ws 87 @(8, 10, 12, 14) @(26000, 26000, 26000, -28496)  # HJLN87
Set-Cells $ws 90 @(8, 10, 12, 14) @(26000, 26000, 78000, -90480)  # HJLN90
Set-Cells $ws 111 @(8, 9, 11, 13) @(4228.5713, 4400, 13200, -10133)  # HIKM111
Set-Cells $ws 137 @(8, 9, 11, 13) @(2469.75, 1917.3529, 5752.0587, -3202.0587)  # HIKM137

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
Set-Cells $ws 32 @(8, 9, 10, 11, 12, 13, 14) @(17980.562, 18952.535, 8584.833000000001, 18952.535, 8584.833000000001, -18665.535, -9158.833000000001)  # HIJKLMN32
Set-Cells $ws 82 @(8, 10, 12, 14) @(40180.832, 40180.832, 40180.832, -40902.832)  # HJLN82
Set-Cells $ws 85 @(8, 10, 12, 14) @(40180.832, 40180.832, 40180.832, -42676.832)  # HJLN85

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
Set-Cells $ws 31 @(8, 9, 10, 11, 12, 13, 14) @(6217.3447, 9845.583000000001, 3656.2354, 9845.583000000001, 3656.2354, -9550.583000000001, -4246.2354)  # HIJKLMN31
Set-Cells $ws 33 @(8, 9, 10, 11, 12, 13) @(3138, 3138, 0, 3138, 0, -2759)  # HIJKLM33
$ws.Cells.Item(33, 14).ClearContents()  # N33
Set-Cells $ws 34 @(8, 9, 10, 11, 12, 13, 14) @(6217.3447, 9845.583000000001, 3656.2354, 9845.583000000001, 3656.2354, -9643.583000000001, -4060.2354)  # HIJKLMN34
Set-Cells $ws 134 @(8, 9, 10, 11, 12, 13, 14) @(2882.6072, 2656.9048, 3559.7144, 7970.714399999999, 10679.1432, -5435.714399999999, -15749.1432)  # HIJKLMN134

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
Set-Cells $ws 88 @(8, 10, 12, 14) @(3273.6875, 3507.4482, 10522.3446, -11378.3446)  # HJLN88
Set-Cells $ws 91 @(8, 10, 12, 14) @(3273.6875, 3507.4482, 10522.3446, -13486.3446)  # HJLN91
Set-Cells $ws 102 @(8, 10, 12, 14) @(5353.231, 5412.727, 16238.181, -21106.181)  # HJLN102
Set-Cells $ws 113 @(8, 9, 10, 11, 12, 13, 14) @(678.25, 691.4681, 609.2222, 2074.4043, 1827.6666, 95.59569999999985, -6167.6666)  # HIJKLMN113
Set-Cells $ws 133 @(8, 9, 10, 11, 12, 13, 14) @(3351.5715, 1871.1818, 4980, 5613.5454, 14940, -553.5454, -25060)  # HIJKLMN133

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
Set-Cells $ws 36 @(8, 9, 10, 11, 12, 13, 14) @(2348.4666, 1162.4, 2941.5, 1162.4, 2941.5, -677.4000000000001, -3911.5)  # HIJKLMN36
Set-Cells $ws 43 @(8, 10, 12) @(1000, 0, 0)  # HJL43
$ws.Cells.Item(43, 14).ClearContents()  # N43
Set-Cells $ws 46 @(8, 9, 10, 11, 12, 14) @(26266.666, 0, 26266.666, 0, 26266.666, -26578.666)  # HIJKLN46
$ws.Cells.Item(46, 13).ClearContents()  # M46
Set-Cells $ws 88 @(8, 10, 12, 14) @(40195, 40195, 40195, -41097)  # HJLN88
Set-Cells $ws 91 @(8, 10, 12, 14) @(40195, 40195, 40195, -43315)  # HJLN91
Set-Cells $ws 113 @(8, 9, 10, 11, 12, 13, 14) @(1980.0714, 2018.4166, 1750, 2018.4166, 1750, 151.5834, -6090)  # HIJKLMN113

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
Set-Cells $ws 22 @(8, 9, 10, 11, 12, 13, 14) @(961.5833, 951.6667, 964.8889, 951.6667, 964.8889, -656.6667, -1554.8889)  # HIJKLMN22
Set-Cells $ws 27 @(8, 9, 10, 11, 12, 13, 14) @(961.5833, 951.6667, 964.8889, 951.6667, 964.8889, -844.6667, -1178.8889)  # HIJKLMN27
Set-Cells $ws 46 @(8, 10, 12) @(860, 900, 900)  # HJL46
